$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Status updates: "?" (Em aberto) -> "!" (Concluida) for several activities ---
$ws.Range("F4").Value = "!"
$ws.Range("F7").Value = "!"
$ws.Range("F9").Value = "!"
$ws.Range("F10").Value = "!"
$ws.Range("F11").Value = "!"

# --- Collaborator reassignment ---
$ws.Range("E10").Value = "Neimar"

# --- Fill in descriptions that were placeholder "||" (ditto marks) with the actual table name ---
$ws.Range("B7").Value = "Tabela tipoUtilitario"
$ws.Range("B9").Value = "Tabela Utilitarios"
$ws.Range("B11").Value = "Tabela Material"
$ws.Range("B13").Value = "Tabela NotaFiscal"
$ws.Range("B15").Value = "Tabela ItensNf"

# --- Update activity description text ---
$ws.Range("B19").Value = "Todas as tabelas e funções "

# --- Match formatting: B7 now takes on the plain "Criar" row look (copy from B6) ---
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122) | Out-Null

# --- B15 takes on the same look as B14 (copy format only) ---
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Update the comment on F3 ---
$comment = $ws.Range("F3").Comment
$newText = "Estatus da atividade:`r`n========================`r`n?  = Em aberto`r`n!   = Concluida`r`n?! = Em execução`r`n!!  = Revisado ok`r`n*  =  Dúvida (Aguarda solução)`r`n"
$comment.Text($newText)
$comment.Shape.TextFrame.Characters().Font.Bold = $true

# --- Update sheet view: scroll so row 4 is at top, select B16 ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B16").Select()
